$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "statut_label" column (B) entirely.
$ws.Columns.Item(2).Delete()

# Remove the "results_1y", "results_3y" and "results" columns
# (which, after the deletion above, sit at I:K) - no longer published.
$ws.Range("I1:K7").Delete()

# Recompute the "statut" / "statut_name" values with the corrected
# status codes and wording.
for ($r = 2; $r -le 7; $r++) {
  $old = $ws.Cells.Item($r, 1).Value()

  $ws.Cells.Item($r, 1).NumberFormat = "@"
  $ws.Cells.Item($r, 2).NumberFormat = "@"

  if ($old -eq "+3") {
    $ws.Cells.Item($r, 1).Value = "2"
    $ws.Cells.Item($r, 2).Value = "2: résultats postés ou publiés entre 12 et 36 mois"
  } else {
    $ws.Cells.Item($r, 1).Value = "4"
    $ws.Cells.Item($r, 2).Value = "4: pas de résultats postés ni publiés"
  }

  $ws.Cells.Item($r, 1).ClearFormats()
  $ws.Cells.Item($r, 2).ClearFormats()
}
